$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New working set of sequence rows: image path, cue word, category
$rows = @(
    @(2, "flower/flower021.jpg", "haken", "flower"),
    @(3, "face/face018.jpg", "tagen", "face"),
    @(4, "flower/flower020.jpg", "drehen", "face"),
    @(5, "flower/flower027.jpg", "fesseln", "flower"),
    @(6, "flower/flower015.jpg", "stechen", "flower"),
    @(7, "face/face010.jpg", "lehnen", "face"),
    @(8, "face/face031.jpg", "segeln", "face"),
    @(9, "face/face005.jpg", "wiegen", "flower"),
    @(10, "flower/flower031.jpg", "liefern", "flower"),
    @(11, "face/face027.jpg", "tauschen", "flower"),
    @(12, "flower/flower011.jpg", "gründen", "flower"),
    @(13, "face/face029.jpg", "stärken", "face"),
    @(14, "flower/flower022.jpg", "füllen", "flower"),
    @(15, "face/face013.jpg", "strahlen", "flower"),
    @(16, "face/face008.jpg", "opfern", "face"),
    @(17, "face/face024.jpg", "jubeln", "face"),
    @(18, "flower/flower000.jpg", "runden", "face"),
    @(19, "face/face026.jpg", "laufen", "face"),
    @(20, "flower/flower007.jpg", "währen", "face"),
    @(21, "face/face019.jpg", "loben", "flower"),
    @(22, "flower/flower025.jpg", "regnen", "face"),
    @(23, "face/face016.jpg", "hupen", "face"),
    @(24, "flower/flower001.jpg", "spielen", "flower"),
    @(25, "face/face030.jpg", "bitten", "face"),
    @(26, "face/face021.jpg", "bleiben", "flower"),
    @(27, "flower/flower002.jpg", "enden", "face"),
    @(28, "flower/flower016.jpg", "sieben", "flower"),
    @(29, "flower/flower028.jpg", "sondern", "face"),
    @(30, "face/face003.jpg", "gelten", "flower"),
    @(31, "flower/flower005.jpg", "dauern", "flower"),
    @(32, "flower/flower014.jpg", "hoffen", "flower"),
    @(33, "face/face025.jpg", "backen", "face")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
